$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -1
$ws.Range("G7").Value = 1
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 1
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 1
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G20").Value = 1
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = 1
$ws.Range("G25").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = -1
$ws.Range("F33").Value = -1
$ws.Range("F34").Value = -1
$ws.Range("F35").Value = -1
$ws.Range("F38").Value = 0
$ws.Range("G53").Value = 1
$ws.Range("G73").Value = 1
$ws.Range("G78").Value = 1
$ws.Range("G80").Value = 1
$ws.Range("G90").Value = 1
$ws.Range("G98").Value = 1
$ws.Range("G121").Value = 1
$ws.Range("G123").Value = 1
$ws.Range("G134").Value = 1
$ws.Range("G137").Value = 1
$ws.Range("G146").Value = 1
$ws.Range("G181").Value = 1
